# Power simulations workbook edit
# Commit: "power results for k=2 adjusted such that d=0.084, as when k=3"
#
# Adds a new column F "new n (based on d=0.084)" with recomputed sample
# sizes (n) for the k=2 rows (rows 2-12), so that they are comparable to
# the k=3 rows which already used d=0.084.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("F1").Value = "new n (based on d=0.084)"

# New n values recomputed for k=2 (rows 2-12) using d = 0.084
$ws.Range("F2").Value  = 10
$ws.Range("F3").Value  = 10
$ws.Range("F4").Value  = 27
$ws.Range("F5").Value  = 55
$ws.Range("F6").Value  = 102
$ws.Range("F7").Value  = 160
$ws.Range("F8").Value  = 257
$ws.Range("F9").Value  = 382
$ws.Range("F10").Value = 577
$ws.Range("F11").Value = 925
$ws.Range("F12").Value = 1892

# Move / update the active cell selection as recorded in the saved file
$ws.Range("I9").Select()
